$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.173.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.789.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.293"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0690"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.782.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.123.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0778"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("E26").Value = "  +1.84%  "
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E33").Value = "  +3.45%  "
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.446.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.652"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("E38").Value = "  +2.58%  "
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.923"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0509"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.49%  "
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("E48").Value = "  -4.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.948.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("E51").Value = "  +0.04%  "
